# "Generate Report for Archive" - refresh the localization-status report:
#  - the zh-cn / de-de status for the sample file flips from
#    "Ready for handoff" to "In Translation"
#  - the Status column on each language sheet (and the corresponding
#    status columns on the Overview sheet) narrows to fit the new,
#    shorter text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update the status text wherever "Ready for handoff" appeared.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# The status columns shrink now that the text is shorter. 12.5 is the
# closest ColumnWidth input that lands on the narrower target width.
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
